$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): wrap existing header text with surrounding spaces ---
$ws.Range("B1").Value = " Task Name "
$ws.Range("C1").Value = " Duration "
$ws.Range("D1").Value = " Dependencies "
$ws.Range("E1").Value = " Status "
$ws.Range("F1").Value = " Resources "

# --- Row 3: Project Initiation ---
$ws.Range("B3").Value = " Project Initiation & Planning "
$ws.Range("C3").Value = " 1 week "
$ws.Range("D3").Value = " None "
$ws.Range("E3").Value = " Not Started "
$ws.Range("F3").Value = " Project Manager "

# --- Row 4: Requirements Gathering ---
$ws.Range("B4").Value = " Requirements Gathering & Analysis "
$ws.Range("C4").Value = " 2 weeks "
$ws.Range("D4").Value = " Project Initiation & Planning "
$ws.Range("E4").Value = " Not Started "
$ws.Range("F4").Value = " Project Manager, Development Team "

# --- Row 5: System Architecture Design (was Database Design) ---
$ws.Range("B5").Value = " System Architecture Design "
$ws.Range("C5").Value = " 2 weeks "
$ws.Range("D5").Value = " Requirements Gathering & Analysis "
$ws.Range("E5").Value = " Not Started "
$ws.Range("F5").Value = " Development Team "

# --- Row 6: UI/UX Design (was API Development) ---
$ws.Range("B6").Value = " UI/UX Design "
$ws.Range("C6").Value = " 3 weeks "
$ws.Range("D6").Value = " System Architecture Design "
$ws.Range("E6").Value = " Not Started "
$ws.Range("F6").Value = " UI/UX Designer "

# --- Row 7: Backend Development (was iOS App Development) ---
$ws.Range("B7").Value = " Backend Development "
$ws.Range("C7").Value = " 12 weeks "
$ws.Range("D7").Value = " System Architecture Design "
$ws.Range("E7").Value = " Not Started "
$ws.Range("F7").Value = " Development Team "

# --- Row 8: iOS App Development (was Android App Development) ---
$ws.Range("B8").Value = " iOS App Development "
$ws.Range("C8").Value = " 10 weeks "
$ws.Range("D8").Value = " UI/UX Design, Backend Development "
$ws.Range("E8").Value = " Not Started "
$ws.Range("F8").Value = " iOS Development Team "

# --- Row 9: Android App Development (was Payment Gateway Integration) ---
$ws.Range("B9").Value = " Android App Development "
$ws.Range("C9").Value = " 10 weeks "
$ws.Range("D9").Value = " UI/UX Design, Backend Development "
$ws.Range("E9").Value = " Not Started "
$ws.Range("F9").Value = " Android Development Team "

# --- Row 10: Payment Gateway Integration (was Testing and Quality Assurance) ---
$ws.Range("B10").Value = " Payment Gateway Integration "
$ws.Range("C10").Value = " 4 weeks "
$ws.Range("D10").Value = " Backend Development "
$ws.Range("E10").Value = " Not Started "
$ws.Range("F10").Value = " Development Team "

# --- Row 11: Quality Assurance & Testing (was User Guide and Help Documentation) ---
$ws.Range("B11").Value = " Quality Assurance & Testing "
$ws.Range("C11").Value = " 6 weeks "
$ws.Range("D11").Value = " iOS App Development, Android App Development, Payment Gateway Integration "
$ws.Range("E11").Value = " Not Started "
$ws.Range("F11").Value = " QA Team "

# --- Row 12: User Documentation & Onboarding Guides (was Deployment and Release Planning) ---
$ws.Range("B12").Value = " User Documentation & Onboarding Guides "
$ws.Range("C12").Value = " 2 weeks "
$ws.Range("D12").Value = " Quality Assurance & Testing "
$ws.Range("E12").Value = " Not Started "
$ws.Range("F12").Value = " Technical Writer "

# --- Row 13: Deployment of Backend (was App Store Submission) ---
$ws.Range("B13").Value = " Deployment of Backend "
$ws.Range("C13").Value = " 1 week "
$ws.Range("D13").Value = " Quality Assurance & Testing "
$ws.Range("E13").Value = " Not Started "
$ws.Range("F13").Value = " DevOps Team "

# --- Row 14: App Store Submission (iOS & Android) (was Marketing and Launch) ---
$ws.Range("B14").Value = " App Store Submission (iOS & Android) "
$ws.Range("C14").Value = " 1 week "
$ws.Range("D14").Value = " Quality Assurance & Testing, Deployment of Backend "
$ws.Range("E14").Value = " Not Started "
$ws.Range("F14").Value = " Project Manager "

# --- Row 15: Marketing & Launch Preparation (was Post-Launch Monitoring & Support) ---
$ws.Range("B15").Value = " Marketing & Launch Preparation "
$ws.Range("C15").Value = " 4 weeks "
$ws.Range("D15").Value = " User Documentation & Onboarding Guides, Deployment of Backend, App Store Submission (iOS & Android) "
$ws.Range("E15").Value = " Not Started "
$ws.Range("F15").Value = " Marketing Team "

# --- Row 16: Project Launch (new row) ---
# Touch A16/G16 so the (already blank) placeholder cells are kept/created
# without assigning an actual empty-string value (which Excel COM ignores
# for previously-unused cells).
$ws.Range("A16").Font.Bold = $false
$ws.Range("B16").Value = " Project Launch "
$ws.Range("C16").Value = " 1 day "
$ws.Range("D16").Value = " Marketing & Launch Preparation "
$ws.Range("E16").Value = " Not Started "
$ws.Range("F16").Value = " Project Manager, Marketing Team "
$ws.Range("G16").Font.Bold = $false

# --- Row 17: Post-Launch Monitoring & Bug Fixes (new row) ---
$ws.Range("A17").Font.Bold = $false
$ws.Range("B17").Value = " Post-Launch Monitoring & Bug Fixes "
$ws.Range("C17").Value = " Ongoing "
$ws.Range("D17").Value = " Project Launch "
$ws.Range("E17").Value = " Not Started "
$ws.Range("F17").Value = " Development Team, QA Team "
$ws.Range("G17").Font.Bold = $false

# --- Row 18: trailing blank marker row (mirrors the old trailing row 16) ---
$ws.Range("A18").Font.Bold = $false

Write-Output "Edit applied"
